$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AA2" = 220
    "AB2" = 60
    "AC2" = 14
    "AD2" = 40
    "F2" = 2.92
    "G2" = 3.35
    "I2" = 3.1
    "J2" = 2.78
    "L2" = 1.6
    "M2" = 1.13
    "N2" = 2.44
    "O2" = 1.59
    "P2" = 1.47
    "Q2" = 2.82
    "R2" = 1.18
    "S2" = 5.9
    "T2" = 2.08
    "U2" = 1.71
    "V2" = 1.51
    "W2" = 1.43
    "X2" = 15
    "Z2" = 55
    "AD3" = 21
    "AF3" = 12
    "AG3" = 11.5
    "AH3" = 46
    "AK3" = 34
    "AN3" = 27
    "F3" = 2.04
    "G3" = 2.16
    "H3" = 4.3
    "J3" = 3.1
    "K3" = 3.5
    "L3" = 1.54
    "N3" = 2.82
    "O3" = 1.49
    "Q3" = 2.44
    "S3" = 5
    "U3" = 1.75
    "W3" = 1.87
    "X3" = 10
    "Z3" = 85
    "AA4" = 34
    "AB4" = 10.5
    "AC4" = 7.2
    "AD4" = 13
    "AE4" = 38
    "AF4" = 28
    "AG4" = 20
    "AH4" = 30
    "AI4" = 200
    "AJ4" = 130
    "AK4" = 95
    "AL4" = 1000
    "AM4" = 290
    "AN4" = 160
    "AO4" = 70
    "F4" = 4.1
    "G4" = 4.7
    "H4" = 2.24
    "I4" = 2.36
    "J4" = 2.88
    "K4" = 3.15
    "M4" = 1.15
    "O4" = 1.65
    "P4" = 1.43
    "T4" = 2.3
    "V4" = 1.74
    "W4" = 1.28
    "X4" = 7.2
    "Y4" = 6.6
    "Z4" = 13.5
    "AA5" = 95
    "AC5" = 15
    "AE5" = 85
    "AF5" = 65
    "AI5" = 200
    "AK5" = 160
    "AL5" = 250
    "F5" = 3.35
    "G5" = 3.95
    "H5" = 2.28
    "J5" = 3
    "K5" = 3.45
    "L5" = 1.55
    "M5" = 1.11
    "N5" = 2.72
    "O5" = 1.48
    "P5" = 1.59
    "Q5" = 2.42
    "S5" = 5
    "T5" = 2.02
    "U5" = 1.8
    "V5" = 1.64
    "W5" = 1.35
    "X5" = 16.5
    "AA6" = 170
    "AB6" = 19
    "AC6" = 14
    "AE6" = 130
    "AI6" = 250
    "AJ6" = 130
    "AK6" = 120
    "AL6" = 260
    "AO6" = 600
    "F6" = 2.6
    "G6" = 2.96
    "I6" = 3.4
    "J6" = 2.92
    "K6" = 3.3
    "L6" = 1.6
    "M6" = 1.12
    "N6" = 2.52
    "O6" = 1.56
    "P6" = 1.5
    "Q6" = 2.64
    "S6" = 5.6
    "T6" = 2.1
    "U6" = 1.75
    "V6" = 1.42
    "W6" = 1.52
    "X6" = 14
    "Y6" = 24
    "AA7" = 70
    "AB7" = 8.6
    "AD7" = 15
    "AE7" = 50
    "AF7" = 14
    "AG7" = 11.5
    "AI7" = 65
    "AJ7" = 36
    "AK7" = 30
    "AL7" = 50
    "AM7" = 140
    "AN7" = 29
    "AO7" = 60
    "F7" = 2.46
    "G7" = 2.5
    "H7" = 3.45
    "I7" = 3.65
    "J7" = 3.15
    "N7" = 3.1
    "P7" = 1.68
    "Q7" = 2.38
    "R7" = 1.25
    "S7" = 4.6
    "U7" = 1.96
    "V7" = 1.38
    "W7" = 1.66
    "X7" = 10.5
    "Y7" = 11.5
    "AA8" = 60
    "AI8" = 290
    "AJ8" = 70
    "AL8" = 130
    "AM8" = 360
    "F8" = 3.15
    "I8" = 3.05
    "K8" = 2.82
    "L8" = 1.83
    "M8" = 1.21
    "O8" = 1.86
    "P8" = 1.35
    "S8" = 9
    "T8" = 2.6
    "U8" = 1.57
    "V8" = 1.48
    "X8" = 6
    "Y8" = 7.2
    "AA9" = 95
    "AB9" = 7
    "AC9" = 7.2
    "AD9" = 18
    "AE9" = 70
    "AF9" = 12.5
    "AH9" = 40
    "AI9" = 240
    "AJ9" = 32
    "AK9" = 34
    "AL9" = 65
    "AM9" = 200
    "AN9" = 34
    "AO9" = 390
    "F9" = 2.32
    "G9" = 2.34
    "H9" = 3.85
    "I9" = 4
    "K9" = 3.25
    "N9" = 2.64
    "O9" = 1.57
    "P9" = 1.53
    "Q9" = 2.76
    "S9" = 5.8
    "T9" = 2.18
    "U9" = 1.76
    "V9" = 1.33
    "W9" = 1.74
    "Y9" = 11
    "Z9" = 26
    "AA10" = 120
    "AB10" = 6.8
    "AC10" = 7.6
    "AD10" = 20
    "AE10" = 85
    "AF10" = 11
    "AG10" = 11.5
    "AH10" = 44
    "AI10" = 260
    "AJ10" = 25
    "AK10" = 65
    "AL10" = 140
    "AM10" = 200
    "AN10" = 25
    "AO10" = 130
    "F10" = 2.06
    "G10" = 2.1
    "H10" = 4.6
    "J10" = 3.25
    "K10" = 3.3
    "M10" = 1.12
    "N10" = 2.88
    "O10" = 1.5
    "P10" = 1.59
    "Q10" = 2.6
    "S10" = 5.2
    "T10" = 2.14
    "U10" = 1.78
    "V10" = 1.26
    "W10" = 1.9
    "Y10" = 13
    "Z10" = 80
    "AA11" = 85
    "AB11" = 7
    "AD11" = 17
    "AE11" = 70
    "AF11" = 13.5
    "AH11" = 38
    "AI11" = 110
    "AJ11" = 36
    "AL11" = 170
    "AM11" = 260
    "AO11" = 110
    "F11" = 2.38
    "G11" = 2.44
    "H11" = 3.65
    "I11" = 3.75
    "J11" = 3.1
    "K11" = 3.2
    "L11" = 1.69
    "M11" = 1.14
    "N11" = 2.46
    "O11" = 1.64
    "P11" = 1.48
    "Q11" = 2.96
    "R11" = 1.17
    "S11" = 6.6
    "T11" = 2.3
    "U11" = 1.72
    "V11" = 1.36
    "W11" = 1.69
    "Y11" = 9.6
    "Z11" = 23
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
